# feat: add 2022-Q4 data
#
# - "总计" sheet: new top row for 2022-Q4 (15 funds, 38.02亿元), with the
#   previously-top 2022-Q3 row pushed down to row 3.
# - New worksheet "2022-Q4" inserted between "总计" and "2022-Q3", holding
#   the per-fund holding detail for the new quarter.
# - Existing "2022-Q3" worksheet (fund holding detail) is left untouched.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the "总计" (totals) sheet.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Copy the existing 2022-Q3 totals row down to row 3 first (values only;
# row 2 will be overwritten with the new 2022-Q4 totals straight after).
$summary.Range("A2").Copy()
$summary.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(3, 2).Value = "2022-Q3"
$summary.Cells.Item(3, 3).Value = 17
$summary.Cells.Item(3, 4).Value = 50.42

# Overwrite row 2 with the new 2022-Q4 totals (index stays 0).
$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 15
$summary.Cells.Item(2, 4).Value = 38.02

# ---------------------------------------------------------------------
# 2) Insert the new "2022-Q4" fund-detail worksheet, right before the
#    existing "2022-Q3" sheet (which is currently the active sheet).
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($q3)
$q4.Name = "2022-Q4"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q4.Cells.Item(1, $i + 2).Value = $headers[$i]
}

$rows = @(
    @{ A=0; B='008903'; C='广发科技先锋混合'; D='129.95'; E='93.81'; F='6.06'; G='7.8750'; H=8 },
    @{ A=1; B='012967'; C='广发行业严选三年持有期混合A'; D='109.75'; E='93.31'; F='5.74'; G='6.2996'; H=9 },
    @{ A=2; B='005911'; C='广发双擎升级混合A'; D='98.42'; E='94.38'; F='5.90'; G='5.8068'; H=8 },
    @{ A=3; B='162703'; C='广发小盘成长混合（LOF）A'; D='96.14'; E='93.72'; F='5.68'; G='5.4608'; H=8 },
    @{ A=4; B='002939'; C='广发创新升级灵活配置混合'; D='78.68'; E='94.25'; F='6.23'; G='4.9018'; H=9 },
    @{ A=5; B='012079'; C='信澳新能源精选混合'; D='40.55'; E='93.85'; F='9.63'; G='3.9050'; H=3 },
    @{ A=6; B='008638'; C='广发科技创新混合A'; D='33.72'; E='91.92'; F='3.71'; G='1.2510'; H=10 },
    @{ A=7; B='004854'; C='广发中证全指汽车指数A'; D='15.08'; E='92.79'; F='5.36'; G='0.8083'; H=6 },
    @{ A=8; B='004855'; C='广发中证全指汽车指数C'; D='11.13'; E='92.79'; F='5.36'; G='0.5966'; H=6 },
    @{ A=9; B='012968'; C='广发行业严选三年持有期混合C'; D='10.36'; E='93.31'; F='5.74'; G='0.5947'; H=9 },
    @{ A=10; B='009132'; C='广发小盘成长混合（LOF）C'; D='3.80'; E='93.72'; F='5.68'; G='0.2158'; H=8 },
    @{ A=11; B='009314'; C='广发双擎升级混合C'; D='3.31'; E='94.38'; F='5.90'; G='0.1953'; H=8 },
    @{ A=12; B='013533'; C='广发科技创新混合C'; D='2.90'; E='91.92'; F='3.71'; G='0.1076'; H=10 },
    @{ A=13; B='011987'; C='财通资管智选核心回报6个月持有期混合A'; D='0.11'; E='38.44'; F='1.51'; G='0.0017'; H=5 },
    @{ A=14; B='011988'; C='财通资管智选核心回报6个月持有期混合C'; D='0.01'; E='38.44'; F='1.51'; G='0.0002'; H=5 }
)

$r = 2
foreach ($row in $rows) {
    $q4.Cells.Item($r, 1).Value = $row.A
    $q4.Cells.Item($r, 2).Value = "'" + $row.B
    $q4.Cells.Item($r, 3).Value = $row.C
    $q4.Cells.Item($r, 4).Value = "'" + $row.D
    $q4.Cells.Item($r, 5).Value = "'" + $row.E
    $q4.Cells.Item($r, 6).Value = "'" + $row.F
    $q4.Cells.Item($r, 7).Value = "'" + $row.G
    $q4.Cells.Item($r, 8).Value = $row.H
    $r = $r + 1
}

# Match the look of the other sheets: bold/boxed header row + index column,
# copied straight from the "总计" sheet's own header styling.
$summary.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$summary.Range("A2").Copy()
$q4.Range("A2:A16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Restore the index values the paste-format pass left untouched (format-only
# paste does not touch cell contents, but re-assert them to be safe).
for ($i = 0; $i -lt $rows.Length; $i++) {
    $q4.Cells.Item($i + 2, 1).Value = $rows[$i].A
}
